{"js": "// First part of 2ch\n//\n// 1) The \"_GoBack\" bookmark that currently sits in the middle of the\n//    sentence \"...terjadi dalam se|buah program...\" (splitting it into\n//    two runs) is removed, and the two runs are merged back into a\n//    single run containing the full, unbroken sentence text.\n// 2) The \"_GoBack\" bookmark is re-created at the end of the \"Studi\n//    Literatur\" list item (after the run, right before the paragraph\n//    mark) - this is where Word's automatic \"last edit\" bookmark now\n//    belongs.\n\n// --- Step 1: remove the old bookmark and merge the split runs -------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst fullSentence =\n  \"yang hanya disadari jauh dalam sebuah proses pengembangan. Salah satu \" +\n  \"metode yang dikembangkan oleh para ilmuwan komputasi untuk bisa \" +\n  \"menangkap sebanyak-banyaknya atau bahkan seluruh kesalahan yang mungkin \" +\n  \"terjadi dalam sebuah program pada tahap yang seawal mungkin adalah \" +\n  \"Metode Formal.\";\n\nconst sentenceResults = context.document.body.search(fullSentence, {\n  matchCase: true,\n  matchWholeWord: false,\n});\nsentenceResults.load(\"items\");\nawait context.sync();\n\n// Re-writing the whole (previously split) sentence as one Replace call\n// makes the host normalise/merge the two adjoining runs that used to be\n// separated by the bookmark.\nsentenceResults.items[0].insertText(fullSentence, \"Replace\");\nawait context.sync();\n\n// --- Step 2: add the bookmark back at the end of \"Studi Literatur\" --------\nconst liResults = context.document.body.search(\"Studi Literatur\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nliResults.load(\"items\");\nawait context.sync();\n\nconst endOfLi = liResults.items[0].getRange(\"End\");\nendOfLi.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# First part of 2ch\n#\n# 1) The \"_GoBack\" bookmark that currently sits in the middle of the\n#    sentence \"...terjadi dalam se|buah program...\" (splitting it into\n#    two runs) is removed, and the two runs are merged back into a\n#    single run containing the full, unbroken sentence text.\n# 2) The \"_GoBack\" bookmark is re-created at the end of the \"Studi\n#    Literatur\" list item (after the run, right before the paragraph\n#    mark) - this is where Word's automatic \"last edit\" bookmark now\n#    belongs.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: remove the old bookmark and merge the split runs -------------\n$bm = $d.Bookmarks(\"_GoBack\")\n$s = $bm.Start\n$bm.Delete()\n\n# Collapse a tiny range around the old bookmark position and run a\n# Find/Replace over it (same text in, same text out). That is enough to\n# make the engine normalise/merge the two adjoining runs that used to be\n# split by the bookmark, while staying narrow enough not to touch the\n# other, earlier occurrence of the same \"dalam se\" substring in this\n# paragraph.\n$r = $d.Range($s - 2, $s + 2)\n$find = $r.Find\n$find.Text = $r.Text\n$find.Replacement.Text = $r.Text\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# --- Step 2: add the bookmark back at the end of \"Studi Literatur\" --------\n$target = $d.Content\n$find2 = $target.Find\n$find2.Text = \"Studi Literatur\"\n$find2.Execute() | Out-Null\n$target.Collapse(0)\n\n# Insert a single, throwaway marker character right at the collapsed\n# insertion point so we get a genuine (non-empty) Range to hand to\n# Bookmarks.Add, then delete the marker again, leaving the bookmark\n# behind in the correct spot.\n$target.InsertAfter(\"@\")\n\n$marker = $d.Content\n$fm = $marker.Find\n$fm.Text = \"@\"\n$fm.Execute() | Out-Null\n\n$d.Bookmarks.Add(\"_GoBack\", $marker)\n$marker.Text = \"\"\n"}
